$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 466, pushing the existing data (and every
# row below it) down by one. Excel copies the formatting of the row above
# onto the freshly-inserted row, which is what we want for column D's
# custom date-number style.
$ws.Rows.Item(466).Insert()

# Populate the new row 466 with the new weekly price-report entry.
$ws.Cells.Item(466, 1).Value = 11
$ws.Cells.Item(466, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(466, 3).Value = "Bíobío"
$ws.Cells.Item(466, 4).Value = 45258
$ws.Cells.Item(466, 5).Value = 8
$ws.Cells.Item(466, 6).Value = "Fruta"
$ws.Cells.Item(466, 7).Value = 100102
$ws.Cells.Item(466, 8).Value = "Cítricos"
$ws.Cells.Item(466, 9).Value = 100102005
$ws.Cells.Item(466, 10).Value = "Naranja"
$ws.Cells.Item(466, 11).Value = "Valencia"
$ws.Cells.Item(466, 12).Value = "Primera"
$ws.Cells.Item(466, 13).Value = 100
$ws.Cells.Item(466, 14).Value = 12000
$ws.Cells.Item(466, 15).Value = 13000
$ws.Cells.Item(466, 16).Value = 12500
$ws.Cells.Item(466, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(466, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(466, 19).Value = 833
$ws.Cells.Item(466, 20).Value = 15
